{"js": "// Add a comma after \"Also\" and append a new closing sentence to the\n// paragraph that ends with \"...next increment quickly.\"\n\nconst body = context.document.body;\n\n// 1) \"Also with our first release\" -> \"Also, with our first release\"\nconst alsoHits = body.search(\"Also with our first release\", { matchCase: true });\nalsoHits.load(\"text\");\nawait context.sync();\n\nif (alsoHits.items.length > 0) {\n  const hit = alsoHits.items[0];\n  const alsoOnly = hit.search(\"Also\", { matchCase: true });\n  alsoOnly.load(\"text\");\n  await context.sync();\n  if (alsoOnly.items.length > 0) {\n    alsoOnly.items[0].insertText(\",\", \"After\");\n    await context.sync();\n  }\n}\n\n// 2) Append new sentence after \"...next increment quickly.\"\nconst endHits = body.search(\"next increment quickly.\", { matchCase: true });\nendHits.load(\"text\");\nawait context.sync();\n\nif (endHits.items.length > 0) {\n  const endHit = endHits.items[endHits.items.length - 1];\n  endHit.insertText(\n    \" These reasons support our choice to use a Incremental Process with the Agile Model.\",\n    \"After\"\n  );\n  await context.sync();\n}\n", "ps1": "# Add a comma after \"Also\" and append a closing sentence to the\n# paragraph ending \"...next increment quickly.\"\n\n$d = $word.ActiveDocument\n\n# 1) \"Also with our first release\" -> \"Also, with our first release\"\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$findRange.Find.Text = \"Also with our first release\"\n$findRange.Find.Replacement.ClearFormatting()\n$findRange.Find.Replacement.Text = \"Also, with our first release\"\n$findRange.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n\n# 2) Append new sentence right after \"...next increment quickly.\"\n$endRange = $d.Content\n$endRange.Find.ClearFormatting()\n$endRange.Find.Replacement.ClearFormatting()\n$endRange.Find.Text = \"next increment quickly.\"\n$endRange.Find.Replacement.Text = \"next increment quickly. These reasons support our choice to use a Incremental Process with the Agile Model.\"\n$endRange.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n"}
